$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'304.47"
$ws.Range("E2").Value = "'0.79%"
$ws.Range("D3").Value = "'35.64"
$ws.Range("E3").Value = "'-4.79%"
$ws.Range("D4").Value = "'5.094"
$ws.Range("E4").Value = "'1.36%"
$ws.Range("D5").Value = "'0.07863"
$ws.Range("E5").Value = "'0.46%"
$ws.Range("D6").Value = "'2.121"
$ws.Range("E6").Value = "'-3.26%"
$ws.Range("D7").Value = "'7.913"
$ws.Range("E7").Value = "'-1.55%"
$ws.Range("D8").Value = "'0.9182"
$ws.Range("E8").Value = "'0.47%"
$ws.Range("D9").Value = "'0.09733"
$ws.Range("E9").Value = "'-0.03%"
$ws.Range("D10").Value = "'0.1856"
$ws.Range("E10").Value = "'-1.95%"
$ws.Range("D11").Value = "'0.08570"
$ws.Range("E11").Value = "'-0.27%"
$ws.Range("D12").Value = "'0.03554"
$ws.Range("E12").Value = "'0.87%"
$ws.Range("D13").Value = "'0.09931"
$ws.Range("E13").Value = "'-0.37%"
$ws.Range("D14").Value = "'0.001443"
$ws.Range("E14").Value = "'-2.66%"
$ws.Range("D15").Value = "'0.005658"
$ws.Range("E15").Value = "'-0.06%"
$ws.Range("D16").Value = "'3.459"
$ws.Range("E16").Value = "'-0.13%"
$ws.Range("D17").Value = "'4.099"
$ws.Range("E17").Value = "'1.58%"
$ws.Range("D18").Value = "'2.553"
$ws.Range("E18").Value = "'23.34%"
$ws.Range("D19").Value = "'0.3426"
$ws.Range("E19").Value = "'-1.07%"
$ws.Range("D20").Value = "'5.217"
$ws.Range("E20").Value = "'9.78%"
$ws.Range("E21").Value = "'0.65%"
$ws.Range("E22").Value = "'-0.04%"
$ws.Range("D23").Value = "'0.04550"
$ws.Range("E23").Value = "'-1.94%"
$ws.Range("E24").Value = "'5.30%"
$ws.Range("D25").Value = "'0.001237"
$ws.Range("E25").Value = "'0.70%"
$ws.Range("D27").Value = "'0.0004753"
$ws.Range("E27").Value = "'0.03%"
$ws.Range("E39").Value = "'4.53%"
$ws.Range("E40").Value = "'-0.09%"
$ws.Range("D41").Value = "'0.007533"
$ws.Range("E41").Value = "'-6.48%"
$ws.Range("D42").Value = "'0.1397"
$ws.Range("E42").Value = "'0.32%"
$ws.Range("D43").Value = "'0.007751"
$ws.Range("E43").Value = "'1.11%"
$ws.Range("D44").Value = "'0.002204"
$ws.Range("E44").Value = "'1.09%"
$ws.Range("E45").Value = "'5.66%"
$ws.Range("D46").Value = "'0.00006320"
$ws.Range("E46").Value = "'4.62%"
$ws.Range("D48").Value = "'0.0005804"
$ws.Range("E48").Value = "'0.06%"
$ws.Range("D49").Value = "'46.59"
$ws.Range("E49").Value = "'623.84%"
$ws.Range("D50").Value = "'0.002001"
$ws.Range("E50").Value = "'-25.62%"
